$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.323.43"
$ws.Range("E2").Value = "  -3.16%  "
$ws.Range("D3").Value = "1.936.84"
$ws.Range("E3").Value = "  -3.21%  "
$ws.Range("D5").Value = "'250.42"
$ws.Range("E5").Value = "  -2.09%  "
$ws.Range("D6").Value = "'0.7229"
$ws.Range("E6").Value = "  -6.79%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.3316"
$ws.Range("E8").Value = "  -4.72%  "
$ws.Range("D9").Value = "'27.68"
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("D10").Value = "'0.07293"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").Value = "'0.8097"
$ws.Range("E11").Value = "  -4.50%  "
$ws.Range("D12").Value = "'0.08094"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").Value = "1.936.21"
$ws.Range("E13").Value = "  -3.25%  "
$ws.Range("D14").Value = "'5.500"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").Value = "'94.67"
$ws.Range("E15").Value = "  -6.40%  "
$ws.Range("D16").Value = "'15.15"
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("D17").Value = "30.332.18"
$ws.Range("E17").Value = "  -3.11%  "
$ws.Range("D18").Value = "'0.000008335"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").Value = "'251.34"
$ws.Range("E19").Value = "  -7.87%  "
$ws.Range("D20").Value = "'5.895"
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").Value = "2.191.48"
$ws.Range("E21").Value = "  -3.14%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'0.9998"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'7.002"
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("D25").Value = "'9.770"
$ws.Range("E25").Value = "  -3.40%  "
$ws.Range("D26").Value = "'163.94"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("D27").Value = "'2.388"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  -3.66%  "
$ws.Range("D29").Value = "'0.1321"
$ws.Range("E29").Value = "  -6.91%  "
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("D31").Value = "'1.348"
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("D32").Value = "'4.433"
$ws.Range("E32").Value = "  -4.82%  "
$ws.Range("D33").Value = "'4.178"
$ws.Range("E33").Value = "  -6.90%  "
$ws.Range("D34").Value = "'0.05207"
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("D35").Value = "'1.288"
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("D36").Value = "'0.7501"
$ws.Range("E36").Value = "  -4.95%  "
$ws.Range("D37").Value = "'2.744"
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").Value = "'0.01977"
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("D39").Value = "'2.825"
$ws.Range("E39").Value = "  -3.74%  "
$ws.Range("D40").Value = "'79.31"
$ws.Range("E40").Value = "  -8.28%  "
$ws.Range("D41").Value = "'6.340"
$ws.Range("E41").Value = "  -7.41%  "
$ws.Range("D42").Value = "'0.4533"
$ws.Range("E42").Value = "  -3.63%  "
$ws.Range("D43").Value = "'2.024"
$ws.Range("E43").Value = "  -6.22%  "
$ws.Range("D44").Value = "'0.8481"
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'101.80"
$ws.Range("E46").Value = "  -3.69%  "
$ws.Range("D47").Value = "'9.722"
$ws.Range("E47").Value = "  -5.45%  "
$ws.Range("D48").Value = "'7.468"
$ws.Range("E48").Value = "  -4.18%  "
$ws.Range("D49").Value = "'36.78"
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("D50").Value = "'0.4189"
$ws.Range("E50").Value = "  -4.18%  "
$ws.Range("D51").Value = "'0.06037"
$ws.Range("E51").Value = "  -0.33%  "